# Refresh the crypto symbol list (GitHub Actions scheduled update).
# Price (col D) and Volume(1h) (col E) are stored as literal text in this
# sheet (e.g. "310.13", "1.23%"), so every new value is written with a
# leading apostrophe to force text entry and avoid Excel's automatic
# number/percentage coercion (which would also lose formats like the
# leading zeros / trailing zeros / thousands separators already present).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'310.22"
$ws.Range("E2").Value = "'1.27%"

# Row 3 - OKB
$ws.Range("E3").Value = "'1.98%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.123"
$ws.Range("E4").Value = "'0.16%"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.07691"
$ws.Range("E5").Value = "'1.34%"

# Row 6 - GateToken
$ws.Range("D6").Value = "'4.284"
$ws.Range("E6").Value = "'0.30%"

# Row 7 - FTXToken
$ws.Range("D7").Value = "'1.624"
$ws.Range("E7").Value = "'0.90%"

# Row 8 - MXToken
$ws.Range("D8").Value = "'0.9212"
$ws.Range("E8").Value = "'1.64%"

# Row 9 - BTSEToken
$ws.Range("E9").Value = "'1.73%"

# Row 10 - LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.1219"
$ws.Range("E10").Value = "'20.72%"

# Row 11 - WazirX
$ws.Range("D11").Value = "'0.1827"
$ws.Range("E11").Value = "'4.27%"

# Row 12 - MandalaExchangeToken
$ws.Range("D12").Value = "'0.09128"
$ws.Range("E12").Value = "'0.52%"

# Row 13 - BitrueCoin
$ws.Range("E13").Value = "'1.08%"

# Row 14 - BitMartToken
$ws.Range("E14").Value = "'-0.41%"

# Row 15 - BitForexToken
$ws.Range("D15").Value = "'0.001245"
$ws.Range("E15").Value = "'-0.04%"

# Row 16 - TigerCash
$ws.Range("D16").Value = "'0.005860"
$ws.Range("E16").Value = "'-0.25%"

# Row 17 - LEO
$ws.Range("E17").Value = "'0.13%"

# Row 19 - MCDex
$ws.Range("D19").Value = "'6.918"
$ws.Range("E19").Value = "'4.15%"

# Row 20 - ProBitToken
$ws.Range("D20").Value = "'0.1387"
$ws.Range("E20").Value = "'2.20%"

# Row 21 - ZBToken
$ws.Range("D21").Value = "'0.2674"
$ws.Range("E21").Value = "'-2.10%"

# Row 22 - CoinExToken
$ws.Range("D22").Value = "'0.04033"
$ws.Range("E22").Value = "'-3.47%"

# Row 23 - BitKan
$ws.Range("D23").Value = "'0.001265"
$ws.Range("E23").Value = "'3.10%"

# Row 24 - HotbitToken
$ws.Range("D24").Value = "'0.004081"
$ws.Range("E24").Value = "'0.68%"

# Row 25 - NitroEx
$ws.Range("D25").Value = "'0.0001266"
$ws.Range("E25").Value = "'-2.78%"

# Row 26 - UpBots
$ws.Range("E26").Value = "'24.53%"

# Row 38 - One
$ws.Range("D38").Value = "'0.02473"
$ws.Range("E38").Value = "'3.70%"

# Row 39 - IDEX
$ws.Range("D39").Value = "'0.05266"
$ws.Range("E39").Value = "'2.49%"

# Row 40 - KickToken
$ws.Range("D40").Value = "'0.007819"
$ws.Range("E40").Value = "'0.50%"

# Row 41 - BKEXToken
$ws.Range("D41").Value = "'0.1312"
$ws.Range("E41").Value = "'1.21%"

# Row 42 - Dexo
$ws.Range("D42").Value = "'0.006784"
$ws.Range("E42").Value = "'-3.79%"

# Row 43 - CEJI
$ws.Range("D43").Value = "'0.001837"
$ws.Range("E43").Value = "'-4.39%"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.008186"
$ws.Range("E44").Value = "'-3.17%"

# Row 45 - PooCoin
$ws.Range("E45").Value = "'-6.38%"

# Row 46 - CoinLion
$ws.Range("D46").Value = "'0.00006791"
$ws.Range("E46").Value = "'6.64%"

# Row 48 - BOLO
$ws.Range("D48").Value = "'0.1803"
$ws.Range("E48").Value = "'2,534.34%"

# Row 49 - CoinbaseStockToken
$ws.Range("D49").Value = "'0.004089"
$ws.Range("E49").Value = "'-7.21%"

# Row 50 - CryptobidCoin
$ws.Range("D50").Value = "'0.00002095"

# Row 51 - SpecialPowerGold
$ws.Range("D51").Value = "'0.0001995"
